$d = $word.ActiveDocument

# 1. Fix typo: "2024game" -> "2024 game"
$d.Paragraphs(1).Range.Text = "logboek: maar 2024 game development"

# 2. Paragraph 4 (empty, right after Week 1 activity line) gains the "tile palette" remark,
#    written as three runs (matching how the author typed it incrementally).
$p = $d.Paragraphs(4)
$p.Range.InsertAfter("Opmerkingen: ik heb een ")
$p = $d.Paragraphs(4)
$p.Range.InsertAfter("tile")
$p = $d.Paragraphs(4)
$p.Range.InsertAfter(" palette gemaakt en een grid toegevoegd voor een mooi level en een achtergrond toegevoegd. ")

# 3. Week 2 section: merge the bullet into the "Activiteiten/Gebeurtenissen:" line and
#    append two new remark paragraphs.
$d.Paragraphs(7).Range.Text = "Activiteiten/Gebeurtenissen: player movement en camera movement gemaakt"
$d.Paragraphs(8).Range.Delete() | Out-Null

$d.Paragraphs(7).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(8).Range.Text = "opmerkingen: ik heb een speler getekend en geanimeerd en daarna ervoor gezorgd dat hij kon lopen en springen. Het viel gelijk op dat de camera niet mee bewoog dus was de rest van de week daarmee bezig."

$d.Paragraphs(8).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(9).Range.Text = "Opmerking mevrouw Jacobs: vond dat het spel geluid mist."

# 4. Week 3 section: the empty paragraph after the activities line gets a remark,
#    followed by two more new paragraphs (a comment, then a blank one).
$d.Paragraphs(13).Range.Text = "Opmerkingen: ik ben bezig geweest met een coin collect system en een death system"

$d.Paragraphs(13).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(14).Range.Text = "Opmerking mevrouw Jacobs: Er mist een Game Over en er mist een finish."

$d.Paragraphs(14).Range.InsertParagraphAfter() | Out-Null

# 5. Week 4 section: append a closing remark paragraph after the activities line.
$d.Paragraphs(17).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(18).Range.Text = "Opmerkingen: Ik heb een game over gemaakt en een health system dat de player vier hartjes heeft en een finish toegevoegd."
